$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel
# auto-converting numeric-looking strings (e.g. "297.03", "-2.05%",
# "12") into real numbers/percentages. We briefly mark the cell as
# Text, assign the literal, then restore the default "Normal" style
# so no visible formatting change is left behind.
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range("D2") "297.03"
Set-TextValue $ws.Range("E2") "-2.05%"
Set-TextValue $ws.Range("G2") "12"

# Row 3
Set-TextValue $ws.Range("D3") "31.24"
Set-TextValue $ws.Range("E3") "-1.68%"
Set-TextValue $ws.Range("G3") "12"

# Row 4
Set-TextValue $ws.Range("D4") "5.106"
Set-TextValue $ws.Range("E4") "-2.38%"
Set-TextValue $ws.Range("G4") "12"

# Row 5
Set-TextValue $ws.Range("D5") "0.07340"
Set-TextValue $ws.Range("E5") "-0.46%"
Set-TextValue $ws.Range("G5") "12"

# Row 6
Set-TextValue $ws.Range("D6") "7.716"
Set-TextValue $ws.Range("E6") "-1.71%"
Set-TextValue $ws.Range("G6") "12"

# Row 7
Set-TextValue $ws.Range("D7") "1.678"
Set-TextValue $ws.Range("E7") "12.62%"
Set-TextValue $ws.Range("G7") "12"

# Row 8
Set-TextValue $ws.Range("D8") "3.727"
Set-TextValue $ws.Range("E8") "-0.14%"
Set-TextValue $ws.Range("G8") "12"

# Row 9
Set-TextValue $ws.Range("D9") "0.9175"
Set-TextValue $ws.Range("E9") "1.06%"
Set-TextValue $ws.Range("G9") "12"

# Row 10
Set-TextValue $ws.Range("E10") "-0.76%"
Set-TextValue $ws.Range("G10") "12"

# Row 11
Set-TextValue $ws.Range("D11") "0.07089"
Set-TextValue $ws.Range("E11") "-5.14%"
Set-TextValue $ws.Range("G11") "12"

# Row 12
Set-TextValue $ws.Range("D12") "0.08061"
Set-TextValue $ws.Range("E12") "1.11%"
Set-TextValue $ws.Range("G12") "12"

# Row 13
Set-TextValue $ws.Range("E13") "0.71%"
Set-TextValue $ws.Range("G13") "12"

# Row 14
Set-TextValue $ws.Range("D14") "0.09899"
Set-TextValue $ws.Range("E14") "0.01%"
Set-TextValue $ws.Range("G14") "12"

# Row 15
Set-TextValue $ws.Range("D15") "0.001496"
Set-TextValue $ws.Range("E15") "0.61%"
Set-TextValue $ws.Range("G15") "12"

# Row 16
Set-TextValue $ws.Range("D16") "0.006240"
Set-TextValue $ws.Range("E16") "-1.88%"
Set-TextValue $ws.Range("G16") "12"

# Row 17
Set-TextValue $ws.Range("D17") "3.449"
Set-TextValue $ws.Range("E17") "-1.13%"
Set-TextValue $ws.Range("G17") "12"

# Row 18
Set-TextValue $ws.Range("E18") "-0.33%"
Set-TextValue $ws.Range("G18") "12"

# Row 19
Set-TextValue $ws.Range("E19") "-2.00%"
Set-TextValue $ws.Range("G19") "12"

# Row 20
Set-TextValue $ws.Range("D20") "0.1348"
Set-TextValue $ws.Range("E20") "1.92%"
Set-TextValue $ws.Range("G20") "12"

# Row 21
Set-TextValue $ws.Range("D21") "4.553"
Set-TextValue $ws.Range("E21") "1.68%"
Set-TextValue $ws.Range("G21") "12"

# Row 22
Set-TextValue $ws.Range("D22") "0.04628"
Set-TextValue $ws.Range("E22") "2.28%"
Set-TextValue $ws.Range("G22") "12"

# Row 23
Set-TextValue $ws.Range("G23") "12"

# Row 24
Set-TextValue $ws.Range("D24") "0.001216"
Set-TextValue $ws.Range("E24") "0.12%"
Set-TextValue $ws.Range("G24") "12"

# Row 25
Set-TextValue $ws.Range("D25") "0.004428"
Set-TextValue $ws.Range("E25") "-0.02%"
Set-TextValue $ws.Range("G25") "12"

# Row 26
Set-TextValue $ws.Range("D26") "0.0001297"
Set-TextValue $ws.Range("G26") "12"

# Row 27
Set-TextValue $ws.Range("E27") "7.80%"
Set-TextValue $ws.Range("G27") "12"

# Row 28
Set-TextValue $ws.Range("G28") "12"

# Row 29
Set-TextValue $ws.Range("G29") "12"

# Row 30
Set-TextValue $ws.Range("G30") "12"

# Row 31
Set-TextValue $ws.Range("G31") "12"

# Row 32
Set-TextValue $ws.Range("G32") "12"

# Row 33
Set-TextValue $ws.Range("G33") "12"

# Row 34
Set-TextValue $ws.Range("G34") "12"

# Row 35
Set-TextValue $ws.Range("G35") "12"

# Row 36
Set-TextValue $ws.Range("G36") "12"

# Row 37
Set-TextValue $ws.Range("G37") "12"

# Row 38
Set-TextValue $ws.Range("G38") "12"

# Row 39
Set-TextValue $ws.Range("D39") "0.01715"
Set-TextValue $ws.Range("E39") "2.62%"
Set-TextValue $ws.Range("G39") "12"

# Row 40
Set-TextValue $ws.Range("D40") "0.04419"
Set-TextValue $ws.Range("E40") "-1.48%"
Set-TextValue $ws.Range("G40") "12"

# Row 41
Set-TextValue $ws.Range("D41") "0.007190"
Set-TextValue $ws.Range("E41") "-0.08%"
Set-TextValue $ws.Range("G41") "12"

# Row 42
Set-TextValue $ws.Range("D42") "0.1329"
Set-TextValue $ws.Range("E42") "-1.10%"
Set-TextValue $ws.Range("G42") "12"

# Row 43
Set-TextValue $ws.Range("D43") "0.002146"
Set-TextValue $ws.Range("E43") "-7.71%"
Set-TextValue $ws.Range("G43") "12"

# Row 44
Set-TextValue $ws.Range("D44") "0.01072"
Set-TextValue $ws.Range("E44") "-24.96%"
Set-TextValue $ws.Range("G44") "12"

# Row 45
Set-TextValue $ws.Range("D45") "0.00006036"
Set-TextValue $ws.Range("E45") "-1.44%"
Set-TextValue $ws.Range("G45") "12"

# Row 46
Set-TextValue $ws.Range("E46") "-21.22%"
Set-TextValue $ws.Range("G46") "12"

# Row 47
Set-TextValue $ws.Range("D47") "1.894"
Set-TextValue $ws.Range("E47") "0.07%"
Set-TextValue $ws.Range("G47") "12"

# Row 48
Set-TextValue $ws.Range("G48") "12"

# Row 49
Set-TextValue $ws.Range("G49") "12"

# Row 50
Set-TextValue $ws.Range("G50") "12"

# Row 51
Set-TextValue $ws.Range("G51") "12"
